$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "'1000"
$ws.Range("B2").Value = "Test Cash"
$ws.Range("C2").Value = "ASSET"
$ws.Range("D2").Value = "Current Asset"
$ws.Range("E2").Value = "'false"
$ws.Range("F2").Value = "'"
$ws.Range("G2").Value = "'true"
$ws.Range("H2").Value = "Cash on hand"
$ws.Range("I2").Value = "'"

# Row 3
$ws.Range("A3").Value = "'1100"
$ws.Range("B3").Value = "Test Bank Account"
$ws.Range("C3").Value = "ASSET"
$ws.Range("D3").Value = "Current Asset"
$ws.Range("E3").Value = "'false"
$ws.Range("F3").Value = "'"
$ws.Range("G3").Value = "'true"
$ws.Range("H3").Value = "Primary bank account"
$ws.Range("I3").Value = "'1000"

# Row 4
$ws.Range("A4").Value = "'1200"
$ws.Range("B4").Value = "Test Accounts Receivable"
$ws.Range("C4").Value = "ASSET"
$ws.Range("D4").Value = "Current Asset"
$ws.Range("E4").Value = "'true"
$ws.Range("F4").Value = "Customer"
$ws.Range("G4").Value = "'true"
$ws.Range("H4").Value = "Amounts owed by customers"
$ws.Range("I4").Value = "'"

# Row 5
$ws.Range("A5").Value = "'2000"
$ws.Range("B5").Value = "Test Accounts Payable"
$ws.Range("C5").Value = "LIABILITY"
$ws.Range("D5").Value = "Current Liability"
$ws.Range("E5").Value = "'true"
$ws.Range("F5").Value = "Vendor"
$ws.Range("G5").Value = "'true"
$ws.Range("H5").Value = "Amounts owed to vendors"
$ws.Range("I5").Value = "'"

# Row 6
$ws.Range("A6").Value = "'2100"
$ws.Range("B6").Value = "Test Credit Card"
$ws.Range("C6").Value = "LIABILITY"
$ws.Range("D6").Value = "Current Liability"
$ws.Range("E6").Value = "'false"
$ws.Range("F6").Value = "'"
$ws.Range("G6").Value = "'true"
$ws.Range("H6").Value = "Business credit card"
$ws.Range("I6").Value = "'2000"

# Row 7
$ws.Range("A7").Value = "'3000"
$ws.Range("B7").Value = "Test Owner Equity"
$ws.Range("C7").Value = "EQUITY"
$ws.Range("D7").Value = "Equity"
$ws.Range("E7").Value = "'false"
$ws.Range("F7").Value = "'"
$ws.Range("G7").Value = "'true"
$ws.Range("H7").Value = "Owner investment"
$ws.Range("I7").Value = "'"

# Row 8
$ws.Range("A8").Value = "'4000"
$ws.Range("B8").Value = "Test Sales Revenue"
$ws.Range("C8").Value = "REVENUE"
$ws.Range("D8").Value = "Operating Revenue"
$ws.Range("E8").Value = "'false"
$ws.Range("F8").Value = "'"
$ws.Range("G8").Value = "'true"
$ws.Range("H8").Value = "Revenue from sales"
$ws.Range("I8").Value = "'"

# Row 9
$ws.Range("A9").Value = "'4100"
$ws.Range("B9").Value = "Test Service Revenue"
$ws.Range("C9").Value = "REVENUE"
$ws.Range("D9").Value = "Operating Revenue"
$ws.Range("E9").Value = "'false"
$ws.Range("F9").Value = "'"
$ws.Range("G9").Value = "'true"
$ws.Range("H9").Value = "Revenue from services"
$ws.Range("I9").Value = "'4000"

# Row 10
$ws.Range("A10").Value = "'5000"
$ws.Range("B10").Value = "Test Rent Expense"
$ws.Range("C10").Value = "EXPENSE"
$ws.Range("D10").Value = "Operating Expense"
$ws.Range("E10").Value = "'false"
$ws.Range("F10").Value = "'"
$ws.Range("G10").Value = "'true"
$ws.Range("H10").Value = "Office rent"
$ws.Range("I10").Value = "'"

# Row 11
$ws.Range("A11").Value = "'5100"
$ws.Range("B11").Value = "Test Utilities Expense"
$ws.Range("C11").Value = "EXPENSE"
$ws.Range("D11").Value = "Operating Expense"
$ws.Range("E11").Value = "'false"
$ws.Range("F11").Value = "'"
$ws.Range("G11").Value = "'true"
$ws.Range("H11").Value = "Electricity, water, etc."
$ws.Range("I11").Value = "'5000"
